$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Header for column C (row 1)
$ws.Range("C1").Value = "ADM1ATIEN"

# Column B: full region names; Column C: title-cased region abbreviations
$ws.Range("B2").Value = "National Capital Region"
$ws.Range("C2").Value = "NCR"

$ws.Range("B3").Value = "Cordillera Administrative Region"
$ws.Range("C3").Value = "CAR"

$ws.Range("B4").Value = "Region I"
$ws.Range("C4").Value = "Ilocos Region"

$ws.Range("B5").Value = "Region II"
$ws.Range("C5").Value = "Cagayan Valley"

$ws.Range("B6").Value = "Region III"
$ws.Range("C6").Value = "Central Luzon"

$ws.Range("B7").Value = "Region IV-A"
$ws.Range("C7").Value = "Calabarzon"

$ws.Range("B8").Value = "Region IV-B"
$ws.Range("C8").Value = "Mimaropa"

$ws.Range("B9").Value = "Region V"
$ws.Range("C9").Value = "Bicol Region"

$ws.Range("B10").Value = "Region VI"
$ws.Range("C10").Value = "Western Visayas"

$ws.Range("B11").Value = "Region VII"
$ws.Range("C11").Value = "Central Visayas"

$ws.Range("B12").Value = "Region VIII"
$ws.Range("C12").Value = "Eastern Visayas"

$ws.Range("B13").Value = "Region IX"
$ws.Range("C13").Value = "Zamboanga Peninsula"

$ws.Range("B14").Value = "Region X"
$ws.Range("C14").Value = "Northern Mindanao"

$ws.Range("B15").Value = "Region XI"
$ws.Range("C15").Value = "Davao Region"

$ws.Range("B16").Value = "Region XII"
$ws.Range("C16").Value = "Soccsksargen"

$ws.Range("B17").Value = "Region XIII"
$ws.Range("C17").Value = "Caraga"

$ws.Range("B18").Value = "Autonomous Region in Muslim Mindanao"
$ws.Range("C18").Value = "ARMM"

# Update the active selection to match the final state of the file
$ws.Range("C7").Select()
